$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the "Tools:" line's two trailing runs into a single run by
# replacing the text spanning both of them (formatting for both runs was
# already identical, so a straight text replace collapses them naturally).
# ---------------------------------------------------------------------------
$tools = $d.Content
$tools.Find.Execute(
    "AWS (RDS, EC2, Aurora), Azure (SQL Db, Managed Instance, VM)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "AWS (RDS, EC2, Aurora), Azure (SQL Db, Managed Instance, VM)", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: rewrite the closing sentence of the "Final Note" paragraph.
# ---------------------------------------------------------------------------
$closing = $d.Content
$found = $closing.Find.Execute(
    ". Ready to bring expertise, innovation, and dedication to your organization.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $sentenceStart = $closing.Start

    # Replace the whole sentence text in-place; this keeps the original run's
    # formatting (rPr) intact on this (currently single) run.
    $closing.Text = ". I help organisations build robust, fully automated infrastructure with PowerShell, Python and Ansible."
    $sentenceEnd = $closing.End

    # Re-locate "Python" inside the freshly written sentence and toggle Bold
    # on/off on that sub-range. This forces the interop layer to split the
    # containing run into three runs (before/“Python”/after) while leaving
    # every character's effective formatting unchanged (Bold ends up back at
    # its original value of false).
    $scope = $d.Range($sentenceStart, $sentenceEnd)
    $scope.Find.Execute("Python", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $scope.Bold = $true
    $scope.Bold = $false
}
